$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the "asignation" typo to "assignation" in the lead-status labels
# (F4 and F13 keep their row position, only the shared-string text changes;
#  F3 and F12 keep the exact same displayed text as before)
$ws.Range("F13").Value = "assignation-lead-status[manage]"
$ws.Range("F4").Value = "assignation-lead-status[lead]"

# Update the active selection from E12 to F9
$ws.Range("F9").Select()
